$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Sheet "Overview": add row 4 for the newly handed-back file
# ---------------------------------------------------------------------------
$wsOverview = $wb.Worksheets.Item("Overview")

$wsOverview.Range("A4").Value = "286c803c-0bf4-4366-ae7a-86341470ef7e.md"
$wsOverview.Range("B4").Value = "e2e\286c803c-0bf4-4366-ae7a-86341470ef7e.md"
$wsOverview.Range("C4").Value = ".md"
$wsOverview.Range("E4").Value = "Handed back: in sync with en-US"
$wsOverview.Range("F4").Value = "Handed back: in sync with en-US"
$wsOverview.Range("G4").Value = "2016-10-26 08:17:58"
$wsOverview.Range("G4").NumberFormat = "yyyy-mm-dd HH:mm:ss"

$wsOverview.Hyperlinks.Add(
    $wsOverview.Range("B4"),
    "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/bc6022a08dcdde835deef71a5010e0de3d6cb3d8/e2e/286c803c-0bf4-4366-ae7a-86341470ef7e.md",
    [System.Reflection.Missing]::Value,
    [System.Reflection.Missing]::Value,
    "e2e\286c803c-0bf4-4366-ae7a-86341470ef7e.md"
) | Out-Null

$wsOverview.ListObjects.Item(1).Resize($wsOverview.Range("A1:G4"))

# ---------------------------------------------------------------------------
# Sheet "zh-cn": add row 4
# ---------------------------------------------------------------------------
$wsZh = $wb.Worksheets.Item("zh-cn")

$wsZh.Range("A4").Value = "286c803c-0bf4-4366-ae7a-86341470ef7e.md"
$wsZh.Range("B4").Value = ".md"
$wsZh.Range("C4").Value = "Handed back: in sync with en-US"
$wsZh.Range("D4").Value = "e2e"
$wsZh.Range("E4").Value = "ht"
$wsZh.Range("F4").Value = "'True"
$wsZh.Range("G4").Value = "286c803c-0bf4-4366-ae7a-86341470ef7e.3018ef7368574ef2c87524820f5abbcfc36ce044.zh-cn.xlf"
$wsZh.Range("H4").Value = "2016-10-26 08:17:47"
$wsZh.Range("H4").NumberFormat = "yyyy-mm-dd HH:mm:ss"
$wsZh.Range("I4").Value = "286c803c-0bf4-4366-ae7a-86341470ef7e.md"
$wsZh.Range("J4").Value = "286c803c-0bf4-4366-ae7a-86341470ef7e.3018ef7368574ef2c87524820f5abbcfc36ce044.zh-cn.xlf"
$wsZh.Range("K4").Value = "2016-10-26 08:18:28"
$wsZh.Range("K4").NumberFormat = "yyyy-mm-dd HH:mm:ss"
$wsZh.Range("L4").Value = "'"
$wsZh.Range("M4").Value = "'True"
$wsZh.Range("N4").Value = "'"
$wsZh.Range("O4").Value = "'False"
$wsZh.Range("P4").Value = "'"

$wsZh.Hyperlinks.Add(
    $wsZh.Range("A4"),
    "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/bc6022a08dcdde835deef71a5010e0de3d6cb3d8/e2e/286c803c-0bf4-4366-ae7a-86341470ef7e.md",
    [System.Reflection.Missing]::Value,
    [System.Reflection.Missing]::Value,
    "286c803c-0bf4-4366-ae7a-86341470ef7e.md"
) | Out-Null

$wsZh.Hyperlinks.Add(
    $wsZh.Range("I4"),
    "https://github.com/OpenLocalizationTestOrg/ol-test0-zhcn/blob/3018ef7368574ef2c87524820f5abbcfc36ce044/e2e/286c803c-0bf4-4366-ae7a-86341470ef7e.md",
    [System.Reflection.Missing]::Value,
    [System.Reflection.Missing]::Value,
    "286c803c-0bf4-4366-ae7a-86341470ef7e.md"
) | Out-Null

$wsZh.ListObjects.Item(1).Resize($wsZh.Range("A1:P4"))

# ---------------------------------------------------------------------------
# Sheet "de-de": add row 4
# ---------------------------------------------------------------------------
$wsDe = $wb.Worksheets.Item("de-de")

$wsDe.Range("A4").Value = "286c803c-0bf4-4366-ae7a-86341470ef7e.md"
$wsDe.Range("B4").Value = ".md"
$wsDe.Range("C4").Value = "Handed back: in sync with en-US"
$wsDe.Range("D4").Value = "e2e"
$wsDe.Range("E4").Value = "ht"
$wsDe.Range("F4").Value = "'True"
$wsDe.Range("G4").Value = "286c803c-0bf4-4366-ae7a-86341470ef7e.3018ef7368574ef2c87524820f5abbcfc36ce044.de-de.xlf"
$wsDe.Range("H4").Value = "2016-10-26 08:17:58"
$wsDe.Range("H4").NumberFormat = "yyyy-mm-dd HH:mm:ss"
$wsDe.Range("I4").Value = "286c803c-0bf4-4366-ae7a-86341470ef7e.md"
$wsDe.Range("J4").Value = "286c803c-0bf4-4366-ae7a-86341470ef7e.3018ef7368574ef2c87524820f5abbcfc36ce044.de-de.xlf"
$wsDe.Range("K4").Value = "2016-10-26 08:18:45"
$wsDe.Range("K4").NumberFormat = "yyyy-mm-dd HH:mm:ss"
$wsDe.Range("L4").Value = "'"
$wsDe.Range("M4").Value = "'True"
$wsDe.Range("N4").Value = "'"
$wsDe.Range("O4").Value = "'False"
$wsDe.Range("P4").Value = "'"

$wsDe.Hyperlinks.Add(
    $wsDe.Range("A4"),
    "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/bc6022a08dcdde835deef71a5010e0de3d6cb3d8/e2e/286c803c-0bf4-4366-ae7a-86341470ef7e.md",
    [System.Reflection.Missing]::Value,
    [System.Reflection.Missing]::Value,
    "286c803c-0bf4-4366-ae7a-86341470ef7e.md"
) | Out-Null

$wsDe.Hyperlinks.Add(
    $wsDe.Range("I4"),
    "https://github.com/OpenLocalizationTestOrg/ol-test0-dede/blob/3018ef7368574ef2c87524820f5abbcfc36ce044/e2e/286c803c-0bf4-4366-ae7a-86341470ef7e.md",
    [System.Reflection.Missing]::Value,
    [System.Reflection.Missing]::Value,
    "286c803c-0bf4-4366-ae7a-86341470ef7e.md"
) | Out-Null

$wsDe.ListObjects.Item(1).Resize($wsDe.Range("A1:P4"))
